$wb = $excel.ActiveWorkbook

# --- Update long space-separated per-play/per-game numeric log strings (YDS, ST sheets) ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "10 2 4 9 20 7 7 1 1 3 16 15 -1 -4 3 -4 1 3 6 1 10 1 1 4 2 8 5 2 6 3 3 8 3 1 12 2 17 3 2 0 26 11 1 4 14 2 6 4 14 4 0 9 -2 2 1 2 3 12 4 0 4 0 1 1 -1 0 1 -4 5 5 8 4 5 3 -1 1 3 11 2 -2 2 4 9 2 -1 2 4 6 6 3 1 3 -4 2 3 12 1 1 11 4 9 5 4 3 1 2 4 13 4 9 0 3 3 0 4 3 7 3 1 0 3 8 1 6 2 1 1 -5 1 0 9 5 -2 1 6 1 1 8 3 7 6 3 1 5 12 3 14 1 -3 3 4 0 1 3 2 2 11 14 4 5 0 4 1 1 1 -2 5 4 2 2 -2 5 0 2 7 8 1 7 9 3 0 0 4 28 15 4 2 5 21 -1 6 16 1 12 14 42 4 1 4 1 0 -5 0 4 -4 2 6 5 3 0 1 -2 8 3 3 -1 1 21 9 -2 3 -1 5 0 12 4 19 2 57 5 13 7 2 1 0 5 4 -4 0 0 2 4 0 3 7 1 8 4 2 4 2 6 11 19 3 2 10 1 2 8 4 5 4 6 8 2 4 3 0 2 -2 2 11 1 0 5 9 3 5 6 0 5 1 13 35 5 9 3 4 4 2 3 5 4 -3 7 3 3 6 -3 8 10 0 7 6 -3 5 5 0 27 6 8 0 0 -3 1 0 1 4 5 3 0 6 0 2 8 1 16 5 0 2 1 2 -2 3"
$ws.Range("B3").Value = "22 6 -8 12 12 1 3 6 2 6 7 10 -1 7 11 6 43 3 23 5 5 3 10 5 6 5 3 8 7 9 11 13 3 20 2 20 25 8 46 9 5 6 4 5 11 4 10 8 20 1 6 15 7 2 5 3 5 7 7 5 24 16 7 11 3 4 5 8 1 15 10 5 12 11 14 19 9 2 -3 11 9 6 24 19 22 7 21 5 6 33 17 2 20 11 22 5 9 4 11 9 24 13 25 17 3 8 9 9 10 8 13 4 13 8 4 4 6 21 8 14 5 2 17 2 8 7 18 21 10 0 7 10 4 6 5 3 5 -2 5 4 33 2 9 5 3 15 7 4 3 7 10 3 7 8 10 14 15 4 63 -1 21 12 15 8 20 14 9 7 5 37 6 6 5 4 3 6 5 17 2 18 7 5 8 -1 5 7 7 2 10 5 9 18 4 34 0 8 6 11 11 7 9 17 5 2 -3 7 5 19 0 5 5 9 4 1 30 8 12 9 -1 5 3 2 9 -1 10 0 4 -1 6 24 -4 0 3 6 16 39 9 10 5 4 -6 9 2 0 -1 17 8 0 11 17 15 7 9 9 25 9 28 23 6 5 27 7 6 24 7 8 14 9 5 6 13 9 9 9 10 8 1 13 7 7 2 19 7 5 20 11 -2 9 9 11 13 4 6 14 8 14 7 2 4 5 15 5 13 5 13 5 37 5 15 14 8 22 16 6 5 8 11 -2"
$ws.Range("C2").Value = "0 11 9 1 3 3 38 6 11 3 3 3 2 10 2 7 3 2 1 7 -1 0 0 6 0 1 1 3 6 7 1 3 9 4 1 4 8 6 5 4 3 8 3 1 -4 5 0 3 4 1 8 -1 3 5 2 -2 2 5 0 31 6 2 0 4 8 8 1 9 0 5 7 8 5 12 3 9 4 1 4 4 10 1 4 -1 5 1 11 9 1 4 10 2 3 2 10 1 4 3 3 24 4 -1 5 7 2 23 8 -3 4 -1 4 4 -2 -2 9 4 0 5 5 11 3 6 0 3 3 3 0 5 48 2 4 -1 3 2 2 1 2 1 0 1 0 9 8 5 2 3 0 2 4 6 6 14 4 3 5 1 7 6 3 4 6 2 4 5 18 4 2 1 -2 6 4 4 3 5 10 0 3 1 1 1 4 0 -2 3 8 4 4 0 5 7 11 3 10 5 1 1 8 6 10 4 3 19 1 12 3 14 4 -2 20 4 7 3 10 4 3 21 4 3 5 2 2 5 3 3 11 2 5 5 0 2 9 3 12 26 1 7 4 11 1 0 6 8 4 2 3 3 1 3 2 4 0 -1 14 2 1 1 12 4 16 1 1 5 13 2 0 6 2 5 10 0 0 9 7 20 5 3 11 3 3 4 11 15 5 3 3 7 2 0 0 0 1 2 4 8 4 0 2 6 0 0 3 3 5 4 5 12 -1 9 1 0 1 7 5 1 4 9 3 5 -4 3 10 -1 2 5 1 8 11 4 8 0 1 3 4 4 2 5 1 11 5 3 1 13 9 1 5 6 3 3 9 4 1 3 2 2 1 3 3 -1 6 16 5 5 5 9 14 -1 1 1 6 8 3 4 -1 10 5 1 0 6 2 5 1 22 0 3 23 6 1 2 1 6 3 0 4 12"
$ws.Range("C3").Value = "23 7 5 29 8 7 23 7 14 18 35 8 79 16 2 15 12 11 18 19 4 5 6 9 13 1 10 50 22 14 3 9 5 22 11 11 0 15 8 0 5 16 11 6 21 37 19 29 3 9 41 7 24 5 18 4 4 21 64 8 8 28 32 6 27 7 37 14 4 17 6 3 6 14 22 2 15 4 9 2 11 6 3 14 9 7 16 8 21 6 19 14 24 5 5 1 9 4 3 7 9 34 8 7 23 11 53 14 2 7 16 29 4 5 1 5 5 2 11 11 27 3 22 10 4 9 2 7 13 22 5 14 59 14 5 6 18 5 7 16 0 15 13 17 19 43 15 6 11 4 9 9 -3 7 17 -2 1 7 20 -4 14 6 3 6 5 14 10 9 6 39 10 6 8 9 7 10 13 9 21 20 3 7 17 6 20 3 5 7 9 37 -1 3 11 13 8 7 12 33 2 17 52 17 16 8 9 0 10 19 11 16 12 7 13 12 7 19 0 18 -2 -2 34 6 7 14 6 2 16 15 11 5 19 10 48 7 16 4 21 14 16 5 14 5 -2 11 3 13 32 8 9 8 6 7 7 4 8 7 9 10 16 5 5 4 6 -1 14 13 2 5 24 10 15 0 3 9 20 19 27 16 5 7 6 7 13 11 26"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value = "59 59 62 64 62 61 67 60 64 60 64 64 66 63 64 66 62 60 66 70 67 61 64 61 65 65"
$ws.Range("B5").Value = "16 19 16 19 41 33 33 14 24 13 17 17 17 17 24 23 14 21 26 25 19 15 21 44 17 31"
$ws.Range("B6").Value = "27 19 0 33 21 24 23 21 26 20 23 28 18 27 26 30 41 29 18 29 38 28 28 13 22 23 35 27 22 30 17"
$ws.Range("D3").Value = "57 38 57 52 55 48 25 38 57 59 46 56 54 56 43 67 47 47 44 49 64 63 64 49 36 39 43 58 52 61 52 39 42 49 46 43 33 47 45 43 61 32 59 70 41 45 43 52 49 62 47 57 54 52 40 55 38"
$ws.Range("D4").Value = "7 0 11 2 0 0 0 0 13 29 0 0 6 0 0 0 2 0 0 0 15 11 0 0 0 0 0 12 13 13 0 0 5 13 0 10 0 1 0 12 0 0 0 0 22 5 13 0 -1 16 2 5 2 11 6 12 0"
$ws.Range("D5").Value = "0 0 0 10 0 11 0 0 7 0 22 14 0 0 11 13 14 0 11 0 48 5 0 0 18 0 0 0 17 0 9 0 0 0 0 0 0 0 -4 0 10 14"

# --- Update numeric summary cells (OFF, DEF, ST, TURNS, PEN sheets) ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 157
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 42
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 23
$ws.Range("L2").Value = 252
$ws.Range("M2").Value = 180
$ws.Range("O2").Value = 31
$ws.Range("P2").Value = 14
$ws.Range("Q2").Value = 439
$ws.Range("C3").Value = 157
$ws.Range("E3").Value = 24
$ws.Range("F3").Value = 83
$ws.Range("G3").Value = 42
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 52
$ws.Range("J3").Value = 39
$ws.Range("N3").Value = 17

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 201
$ws.Range("F2").Value = 64
$ws.Range("G2").Value = 53
$ws.Range("J2").Value = 25
$ws.Range("L2").Value = 220
$ws.Range("M2").Value = 140
$ws.Range("O2").Value = 9
$ws.Range("Q2").Value = 424
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 127
$ws.Range("E3").Value = 23
$ws.Range("F3").Value = 93
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = 28
$ws.Range("I3").Value = 47
$ws.Range("J3").Value = 37
$ws.Range("N3").Value = 13

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 50
$ws.Range("D2").Value = 57
$ws.Range("F2").Value = 126
$ws.Range("G2").Value = 118
$ws.Range("J2").Value = 56
$ws.Range("K2").Value = 52
$ws.Range("L2").Value = 27
$ws.Range("M2").Value = 18
$ws.Range("B3").Value = 24

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 6
$ws.Range("D3").Value = 6

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 18
$ws.Range("B3").Value = 16
$ws.Range("D4").Value = 11
